$wb = $excel.ActiveWorkbook

# Insert a new "Player Info" sheet before the first existing sheet
# ("ODI Batting") so the final tab order is:
#   Player Info, ODI Batting, ODI Bowling
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row - keep ID as text to match source formatting
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5953"
$playerInfo.Range("B2").Value = "Shahidullah Kamal"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# Update "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, link -> bare match code
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4525"

# Update "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, link -> bare match code
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4525"
